{"js": "// Rewrites the Bescheid letter body to the \"Level 2\" draft:\n// - adds an \"Einleitung:\" / \"Tenor:\" / \"Unterschrift mit Gru\u00dfformel:\" scaffold\n// - rewords several sentences (subject line, salutation, tenor items)\n// - merges / rewrites the \"Begr\u00fcndung\" paragraphs into the new legal reasoning\n// - rewrites the Rechtsbehelfsbelehrung and replaces \"[Unterschrift]\" with a signature block\n// The whole letter lives in a single paragraph whose run text is interspersed\n// with manual line breaks (<w:br/>), so we rebuild that paragraph's text from\n// scratch as a list of (text, numberOfTrailingBreaks) segments.\n\nconst segments = [\n  [\"Einleitung:\", 2],\n  [\"Landratsamt Ortenaukreis\", 1],\n  [\"Herrn Franz Konrad\", 1],\n  [\"Sachbearbeiter\", 2],\n  [\"Reparatur Ihres Fachwerkhauses in Neuried, Lange Stra\u00dfe 12\", 2],\n  [\"Sehr geehrter Herr Konrad,\", 2],\n  [\"hiermit ergeht folgender Bescheid:\", 2],\n  [\"Tenor:\", 2],\n  [\"1. Sie sind verpflichtet, das Fachwerkhausdach mit Biberschwanz-Dachziegeln zu reparieren.\", 1],\n  [\"2. Die Reparaturanordnung ist sofort vollziehbar.\", 2],\n  [\"Begr\u00fcndung:\", 2],\n  [\"Sie sind als Forstrat zusammen mit Ihrem Bruder, dem Studenten Georg Konrad, Eigent\u00fcmer des o.g. Fachwerkhauses. Das Haus stammt aus dem Jahre 1865 und geh\u00f6rt zu den wenigen voll erhaltenen Exemplaren seiner Art am Oberrhein. Durch einen Sturm wurden ca. 50 Biberschwanz-Dachziegel abgedeckt, wodurch das Fachwerkhausdach besch\u00e4digt wurde.\", 2],\n  [\"Das Fachwerkhaus ist ein Kulturdenkmal, da es nach \u00a7 2 Abs. 1 DSchG ein \u00f6ffentliches Erhaltungsinteresse aus heimatgeschichtlichen Gr\u00fcnden gibt. Das Kulturdenkmal ist gef\u00e4hrdet, da durch das beeintr\u00e4chtigte Erscheinungsbild bereits ein Schaden entstanden ist.\", 2],\n  [\"Die Reparaturanordnung st\u00fctzt sich auf \u00a7 1 Abs. 1 in Verbindung mit \u00a7 7 Abs. 1 DSchG. Danach k\u00f6nnen wir Ma\u00dfnahmen anordnen, wenn ein Kulturdenkmal gef\u00e4hrdet ist. Die Anordnung der BSD ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da eine kosteng\u00fcnstigere Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erf\u00fcllen. Der Vorteil f\u00fcr die Allgemeinheit durch die Ansehnlichkeit des Denkmals rechtfertigt den finanziellen Nachteil f\u00fcr den Eigent\u00fcmer (F.K.).\", 2],\n  [\"Als Pflichtige kommen sowohl F.K als auch G.K in Betracht. F.K k\u00f6nnte pflichtig sein sinngem\u00e4\u00df \u00a7 7 Abs. 1 Satz 1 DSchG und \u00a7 7 PolG, da er Eigent\u00fcmer einer Sache ist, von deren Zustand eine Gefahr ausgeht. Der G.K ist ebenfalls Eigent\u00fcmer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.\", 2],\n  [\"Es besteht keine Problematik im Sinne des \u00a7 21 LVwVfg. Nach \u00a73 Abs. 4 DSchG muss das Landesamt f\u00fcr Denkmalpflege angeh\u00f6rt werden. Nach \u00a7 28 Abs. 1 LVwVfG ist F.K und G.K die Gelegenheit zur \u00c4u\u00dferung zu geben.\", 2],\n  [\"Rechtsbehelfsbelehrung:\", 2],\n  [\"Gegen die Dachdeckungsanordnung k\u00f6nnen Sie innerhalb eines Monats nach Bekanntgabe bei dem Landratsamt Ortenaukreis, Badstra\u00dfe 22, 77652 Offenburg Widerspruch einlegen. Gegen die Anordnung der sofortigen Vollziehung k\u00f6nnen Sie den Antrag auf Wiederherstellung der aufschiebenden Wirkung beim Verwaltungsgericht Freiburg, Herbstburgerstra\u00dfe 115, 79104 Freiburg stellen.\", 2],\n  [\"Unterschrift mit Gru\u00dfformel:\", 2],\n  [\"Mit freundlichen Gr\u00fc\u00dfen\", 2],\n  [\"Landratsamt Ortenaukreis\", 1],\n  [\"Sachbearbeiter\", 0]\n];\n\n// \"\\v\" (vertical tab) is how Office.js represents a manual line break\n// inside a text string; Word turns each one into a <w:br/> on insert.\nlet newText = \"\";\nfor (const [text, breaks] of segments) {\n  newText += text + \"\\v\".repeat(breaks);\n}\n\nconst body = context.document.body;\n\n// The source document is a single paragraph containing a single run full of\n// <w:t>/<w:br/> pairs. Load paragraphs so we can reuse the first paragraph's\n// formatting instead of just appending to a cleared body (clear() keeps the\n// first paragraph's mark/properties intact).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  const firstPara = paragraphs.items[0];\n  firstPara.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n\n  // Remove any other leftover paragraphs (normally there are none, but guard\n  // against the body having more than the single paragraph we expect).\n  paragraphs.load(\"items\");\n  await context.sync();\n  for (let i = paragraphs.items.length - 1; i >= 1; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n} else {\n  body.clear();\n  body.insertText(newText, Word.InsertLocation.start);\n  await context.sync();\n}\n", "ps1": "# Rewrites the Bescheid letter body to the \"Level 2\" draft:\n# - adds an \"Einleitung:\" / \"Tenor:\" / \"Unterschrift mit Gru\u00dfformel:\" scaffold\n# - rewords several sentences (subject line, salutation, tenor items)\n# - merges / rewrites the \"Begruendung\" paragraphs into the new legal reasoning\n# - rewrites the Rechtsbehelfsbelehrung and replaces \"[Unterschrift]\" with a signature block\n#\n# The whole letter lives in a single paragraph whose run text is interspersed\n# with manual line breaks (Word's vertical-tab \"`v\" char, i.e. <w:br/> in the\n# OOXML), so we rebuild the document's text from scratch as a list of\n# (text, numberOfTrailingBreaks) segments and assign it back to Content.Text\n# (this keeps the single paragraph / sectPr exactly as before).\n\n$segments = @(\n    @('Einleitung:', 2),\n    @('Landratsamt Ortenaukreis', 1),\n    @('Herrn Franz Konrad', 1),\n    @('Sachbearbeiter', 2),\n    @('Reparatur Ihres Fachwerkhauses in Neuried, Lange Stra\u00dfe 12', 2),\n    @('Sehr geehrter Herr Konrad,', 2),\n    @('hiermit ergeht folgender Bescheid:', 2),\n    @('Tenor:', 2),\n    @('1. Sie sind verpflichtet, das Fachwerkhausdach mit Biberschwanz-Dachziegeln zu reparieren.', 1),\n    @('2. Die Reparaturanordnung ist sofort vollziehbar.', 2),\n    @('Begr\u00fcndung:', 2),\n    @('Sie sind als Forstrat zusammen mit Ihrem Bruder, dem Studenten Georg Konrad, Eigent\u00fcmer des o.g. Fachwerkhauses. Das Haus stammt aus dem Jahre 1865 und geh\u00f6rt zu den wenigen voll erhaltenen Exemplaren seiner Art am Oberrhein. Durch einen Sturm wurden ca. 50 Biberschwanz-Dachziegel abgedeckt, wodurch das Fachwerkhausdach besch\u00e4digt wurde.', 2),\n    @('Das Fachwerkhaus ist ein Kulturdenkmal, da es nach \u00a7 2 Abs. 1 DSchG ein \u00f6ffentliches Erhaltungsinteresse aus heimatgeschichtlichen Gr\u00fcnden gibt. Das Kulturdenkmal ist gef\u00e4hrdet, da durch das beeintr\u00e4chtigte Erscheinungsbild bereits ein Schaden entstanden ist.', 2),\n    @('Die Reparaturanordnung st\u00fctzt sich auf \u00a7 1 Abs. 1 in Verbindung mit \u00a7 7 Abs. 1 DSchG. Danach k\u00f6nnen wir Ma\u00dfnahmen anordnen, wenn ein Kulturdenkmal gef\u00e4hrdet ist. Die Anordnung der BSD ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da eine kosteng\u00fcnstigere Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erf\u00fcllen. Der Vorteil f\u00fcr die Allgemeinheit durch die Ansehnlichkeit des Denkmals rechtfertigt den finanziellen Nachteil f\u00fcr den Eigent\u00fcmer (F.K.).', 2),\n    @('Als Pflichtige kommen sowohl F.K als auch G.K in Betracht. F.K k\u00f6nnte pflichtig sein sinngem\u00e4\u00df \u00a7 7 Abs. 1 Satz 1 DSchG und \u00a7 7 PolG, da er Eigent\u00fcmer einer Sache ist, von deren Zustand eine Gefahr ausgeht. Der G.K ist ebenfalls Eigent\u00fcmer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.', 2),\n    @('Es besteht keine Problematik im Sinne des \u00a7 21 LVwVfg. Nach \u00a73 Abs. 4 DSchG muss das Landesamt f\u00fcr Denkmalpflege angeh\u00f6rt werden. Nach \u00a7 28 Abs. 1 LVwVfG ist F.K und G.K die Gelegenheit zur \u00c4u\u00dferung zu geben.', 2),\n    @('Rechtsbehelfsbelehrung:', 2),\n    @('Gegen die Dachdeckungsanordnung k\u00f6nnen Sie innerhalb eines Monats nach Bekanntgabe bei dem Landratsamt Ortenaukreis, Badstra\u00dfe 22, 77652 Offenburg Widerspruch einlegen. Gegen die Anordnung der sofortigen Vollziehung k\u00f6nnen Sie den Antrag auf Wiederherstellung der aufschiebenden Wirkung beim Verwaltungsgericht Freiburg, Herbstburgerstra\u00dfe 115, 79104 Freiburg stellen.', 2),\n    @('Unterschrift mit Gru\u00dfformel:', 2),\n    @('Mit freundlichen Gr\u00fc\u00dfen', 2),\n    @('Landratsamt Ortenaukreis', 1),\n    @('Sachbearbeiter', 0)\n)\n\n$newText = ''\nforeach ($seg in $segments) {\n    $segText = $seg[0]\n    $breakCount = $seg[1]\n    $newText += $segText\n    for ($i = 0; $i -lt $breakCount; $i++) {\n        $newText += \"`v\"\n    }\n}\n\n$d = $word.ActiveDocument\n$d.Content.Text = $newText\n"}
